$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2023" column (K) to the right of the existing "2022" column (J),
# copying the formatting/style of column J so the new column matches the
# existing table look (header style, currency/number style, borders, etc.).
$ws.Range("J3:J6").Copy($ws.Range("K3:K6"))

# Now overwrite the copied values with the new 2023 figures.
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 2255.6
$ws.Range("K5").Value = 1202
$ws.Range("K6").Value = 2519.3000000000002
